# Add a new "images" column to the MSDS sheet, listing the picture file
# name that documents the hazard pictogram for each product row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The long free-text cells in row 4 (columns I:Q) were wrapping, which is
# what forced the huge 409.5pt row height. Turn wrapping off and let the
# row shrink back down to its natural height before we touch anything else.
$ws.Range("I4:Q4").WrapText = $false
$ws.Rows.Item(4).AutoFit()

# Insert a new column at D, pushing the existing D:S columns to E:T.
$ws.Range("D1").EntireColumn.Insert()

# Header for the new column, and the image file referenced on the data row.
$ws.Range("D1").Value = "images"
$ws.Range("D4").Value = "fire.png"

# Leave the selection where the author last left it.
$ws.Range("D5").Select() | Out-Null
